$wb = $excel.ActiveWorkbook

# --- NamedThing sheet: drop the "name" and "description" columns, keep only "id" ---
$wsNamedThing = $wb.Worksheets.Item("NamedThing")
$wsNamedThing.Range("B1:C1").EntireColumn.Delete()

# --- BioSample sheet: remove the vital_status data validation, drop the trailing
#     "description" column, and replace the header row with the new field names ---
$wsBioSample = $wb.Worksheets.Item("BioSample")
$wsBioSample.Range("D2:D1048576").Validation.Delete()
$wsBioSample.Range("G1").EntireColumn.Delete()

$wsBioSample.Range("A1").Value = "depth"
$wsBioSample.Range("B1").Value = "sample_type"
$wsBioSample.Range("C1").Value = "latitude"
$wsBioSample.Range("D1").Value = "longitude"
$wsBioSample.Range("E1").Value = "bacteria"
$wsBioSample.Range("F1").Value = "id"

Write-Host "Applied header/column changes to NamedThing and BioSample sheets"
